# Línea 141 schedule refresh - new scrape at 20:12:55 (470/46/59 filas)
# Updates header metadata (last-updated timestamp + row counts) and
# rewrites every schedule row that changed position/value across the
# three sheets (LP1912, LP1912-215, 6203-6173), including newly scraped rows.
$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws = $wb.Worksheets.Item('LP1912')
$ws.Range('A2').Value = 'Última actualización: 20:12:55'
$ws.Range('A3').Value = 'Total filas: 470'

$arr = New-Object 'object[,]' 3,5
$arr[0,0] = '05:49:10'
$arr[0,1] = '07:32'
$arr[0,2] = '84_COLONIA URQUIZA-ESC 49'
$arr[0,3] = 103
$arr[0,4] = 'LP1912'
$arr[1,0] = '05:49:10'
$arr[1,1] = '07:32'
$arr[1,2] = '16_SANTA ANA'
$arr[1,3] = 103
$arr[1,4] = 'LP1912'
$arr[2,0] = '05:49:10'
$arr[2,1] = '07:32'
$arr[2,2] = '11_ETCHEVERRY'
$arr[2,3] = 103
$arr[2,4] = 'LP1912'
$ws.Range('A47:E49').Value = $arr

$arr = New-Object 'object[,]' 2,5
$arr[0,0] = '07:58:19'
$arr[0,1] = '08:42'
$arr[0,2] = '16_SANTA ANA'
$arr[0,3] = 44
$arr[0,4] = 'LP1912'
$arr[1,0] = '06:57:11'
$arr[1,1] = '08:42'
$arr[1,2] = '81_EL PELIGRO'
$arr[1,3] = 105
$arr[1,4] = 'LP1912'
$ws.Range('A78:E79').Value = $arr

$arr = New-Object 'object[,]' 2,5
$arr[0,0] = '14:00:52'
$arr[0,1] = '14:00'
$arr[0,2] = '16_SANTA ANA'
$arr[0,3] = 0
$arr[0,4] = 'LP1912'
$arr[1,0] = '14:00:52'
$arr[1,1] = '14:00'
$arr[1,2] = '14_ABASTO'
$arr[1,3] = 0
$arr[1,4] = 'LP1912'
$ws.Range('A240:E241').Value = $arr

$arr = New-Object 'object[,]' 2,5
$arr[0,0] = '14:00:52'
$arr[0,1] = '15:56'
$arr[0,2] = '27_EL RETIRO'
$arr[0,3] = 116
$arr[0,4] = 'LP1912'
$arr[1,0] = '14:44:25'
$arr[1,1] = '15:56'
$arr[1,2] = '17_ROMERO'
$arr[1,3] = 72
$arr[1,4] = 'LP1912'
$ws.Range('A284:E285').Value = $arr

$arr = New-Object 'object[,]' 2,5
$arr[0,0] = '17:15:09'
$arr[0,1] = '17:16'
$arr[0,2] = '15_ABASTO'
$arr[0,3] = 1
$arr[0,4] = 'LP1912'
$arr[1,0] = '16:52:47'
$arr[1,1] = '17:16'
$arr[1,2] = '11_ETCHEVERRY'
$arr[1,3] = 24
$arr[1,4] = 'LP1912'
$ws.Range('A327:E328').Value = $arr

$arr = New-Object 'object[,]' 2,5
$arr[0,0] = '15:51:48'
$arr[0,1] = '17:38'
$arr[0,2] = '17_ROMERO'
$arr[0,3] = 107
$arr[0,4] = 'LP1912'
$arr[1,0] = '16:18:00'
$arr[1,1] = '17:38'
$arr[1,2] = '27_EL RETIRO'
$arr[1,3] = 80
$arr[1,4] = 'LP1912'
$ws.Range('A342:E343').Value = $arr

$arr = New-Object 'object[,]' 2,5
$arr[0,0] = '17:57:54'
$arr[0,1] = '17:57'
$arr[0,2] = '16_SANTA ANA'
$arr[0,3] = 0
$arr[0,4] = 'LP1912'
$arr[1,0] = '17:57:54'
$arr[1,1] = '17:57'
$arr[1,2] = '11_ETCHEVERRY'
$arr[1,3] = 0
$arr[1,4] = 'LP1912'
$ws.Range('A355:E356').Value = $arr

$arr = New-Object 'object[,]' 2,5
$arr[0,0] = '17:42:01'
$arr[0,1] = '19:30'
$arr[0,2] = '225_GOMEZ'
$arr[0,3] = 108
$arr[0,4] = 'LP1912'
$arr[1,0] = '19:15:23'
$arr[1,1] = '19:30'
$arr[1,2] = '16_SANTA ANA'
$arr[1,3] = 15
$arr[1,4] = 'LP1912'
$ws.Range('A416:E417').Value = $arr

$arr = New-Object 'object[,]' 33,5
$arr[0,0] = '20:12:55'
$arr[0,1] = '20:21'
$arr[0,2] = '11_ETCHEVERRY'
$arr[0,3] = 9
$arr[0,4] = 'LP1912'
$arr[1,0] = '18:37:39'
$arr[1,1] = '20:22'
$arr[1,2] = '11_ETCHEVERRY'
$arr[1,3] = 105
$arr[1,4] = 'LP1912'
$arr[2,0] = '19:39:27'
$arr[2,1] = '20:22'
$arr[2,2] = '16_SANTA ANA'
$arr[2,3] = 43
$arr[2,4] = 'LP1912'
$arr[3,0] = '18:37:39'
$arr[3,1] = '20:23'
$arr[3,2] = '215A_EL PATO'
$arr[3,3] = 106
$arr[3,4] = 'LP1912'
$arr[4,0] = '18:49:07'
$arr[4,1] = '20:24'
$arr[4,2] = '215A_EL PATO'
$arr[4,3] = 95
$arr[4,4] = 'LP1912'
$arr[5,0] = '19:52:49'
$arr[5,1] = '20:30'
$arr[5,2] = '225_GOMEZ'
$arr[5,3] = 38
$arr[5,4] = 'LP1912'
$arr[6,0] = '18:37:39'
$arr[6,1] = '20:31'
$arr[6,2] = '225_GOMEZ'
$arr[6,3] = 114
$arr[6,4] = 'LP1912'
$arr[7,0] = '20:12:55'
$arr[7,1] = '20:34'
$arr[7,2] = '16_SANTA ANA'
$arr[7,3] = 22
$arr[7,4] = 'LP1912'
$arr[8,0] = '20:12:55'
$arr[8,1] = '20:43'
$arr[8,2] = '11_ETCHEVERRY'
$arr[8,3] = 31
$arr[8,4] = 'LP1912'
$arr[9,0] = '19:15:23'
$arr[9,1] = '20:44'
$arr[9,2] = '11_ETCHEVERRY'
$arr[9,3] = 89
$arr[9,4] = 'LP1912'
$arr[10,0] = '20:12:55'
$arr[10,1] = '20:46'
$arr[10,2] = '16_SANTA ANA'
$arr[10,3] = 34
$arr[10,4] = 'LP1912'
$arr[11,0] = '19:39:27'
$arr[11,1] = '20:52'
$arr[11,2] = '23_HERNANDEZ'
$arr[11,3] = 73
$arr[11,4] = 'LP1912'
$arr[12,0] = '19:39:27'
$arr[12,1] = '20:52'
$arr[12,2] = '15_ABASTO'
$arr[12,3] = 73
$arr[12,4] = 'LP1912'
$arr[13,0] = '18:56:08'
$arr[13,1] = '20:53'
$arr[13,2] = '11_ETCHEVERRY'
$arr[13,3] = 117
$arr[13,4] = 'LP1912'
$arr[14,0] = '20:12:55'
$arr[14,1] = '20:55'
$arr[14,2] = '10_OLMOS'
$arr[14,3] = 43
$arr[14,4] = 'LP1912'
$arr[15,0] = '19:39:27'
$arr[15,1] = '20:56'
$arr[15,2] = '27_EL RETIRO'
$arr[15,3] = 77
$arr[15,4] = 'LP1912'
$arr[16,0] = '19:15:23'
$arr[16,1] = '20:56'
$arr[16,2] = '10_OLMOS'
$arr[16,3] = 101
$arr[16,4] = 'LP1912'
$arr[17,0] = '19:15:23'
$arr[17,1] = '20:57'
$arr[17,2] = '27_EL RETIRO'
$arr[17,3] = 102
$arr[17,4] = 'LP1912'
$arr[18,0] = '19:15:23'
$arr[18,1] = '21:04'
$arr[18,2] = '84_COLONIA URQUIZA-ESC 49'
$arr[18,3] = 109
$arr[18,4] = 'LP1912'
$arr[19,0] = '19:52:49'
$arr[19,1] = '21:07'
$arr[19,2] = '215B_EL PATO'
$arr[19,3] = 75
$arr[19,4] = 'LP1912'
$arr[20,0] = '19:15:23'
$arr[20,1] = '21:08'
$arr[20,2] = '215B_EL PATO'
$arr[20,3] = 113
$arr[20,4] = 'LP1912'
$arr[21,0] = '19:52:49'
$arr[21,1] = '21:20'
$arr[21,2] = '26_HERNANDEZ'
$arr[21,3] = 88
$arr[21,4] = 'LP1912'
$arr[22,0] = '19:39:27'
$arr[22,1] = '21:21'
$arr[22,2] = '26_HERNANDEZ'
$arr[22,3] = 102
$arr[22,4] = 'LP1912'
$arr[23,0] = '19:52:49'
$arr[23,1] = '21:22'
$arr[23,2] = '10_OLMOS'
$arr[23,3] = 90
$arr[23,4] = 'LP1912'
$arr[24,0] = '20:12:55'
$arr[24,1] = '21:22'
$arr[24,2] = '15_ABASTO'
$arr[24,3] = 70
$arr[24,4] = 'LP1912'
$arr[25,0] = '19:39:27'
$arr[25,1] = '21:23'
$arr[25,2] = '10_OLMOS'
$arr[25,3] = 104
$arr[25,4] = 'LP1912'
$arr[26,0] = '20:12:55'
$arr[26,1] = '21:37'
$arr[26,2] = '14_ABASTO'
$arr[26,3] = 85
$arr[26,4] = 'LP1912'
$arr[27,0] = '19:52:49'
$arr[27,1] = '21:37'
$arr[27,2] = '17_ROMERO'
$arr[27,3] = 105
$arr[27,4] = 'LP1912'
$arr[28,0] = '19:52:49'
$arr[28,1] = '21:38'
$arr[28,2] = '14_ABASTO'
$arr[28,3] = 106
$arr[28,4] = 'LP1912'
$arr[29,0] = '19:39:27'
$arr[29,1] = '21:38'
$arr[29,2] = '17_ROMERO'
$arr[29,3] = 119
$arr[29,4] = 'LP1912'
$arr[30,0] = '20:12:55'
$arr[30,1] = '21:46'
$arr[30,2] = '215A_EL PATO'
$arr[30,3] = 94
$arr[30,4] = 'LP1912'
$arr[31,0] = '19:52:49'
$arr[31,1] = '21:47'
$arr[31,2] = '215A_EL PATO'
$arr[31,3] = 115
$arr[31,4] = 'LP1912'
$arr[32,0] = '20:12:55'
$arr[32,1] = '22:07'
$arr[32,2] = '17_ROMERO'
$arr[32,3] = 115
$arr[32,4] = 'LP1912'
$ws.Range('A443:E475').Value = $arr

# ---- Sheet: LP1912-215 ----
$ws = $wb.Worksheets.Item('LP1912-215')
$ws.Range('A2').Value = 'Última actualización: 20:12:55'
$ws.Range('A3').Value = 'Total filas: 46'

$arr = New-Object 'object[,]' 2,5
$arr[0,0] = '20:12:55'
$arr[0,1] = '21:46'
$arr[0,2] = '215A_EL PATO'
$arr[0,3] = 94
$arr[0,4] = 'LP1912'
$arr[1,0] = '19:52:49'
$arr[1,1] = '21:47'
$arr[1,2] = '215A_EL PATO'
$arr[1,3] = 115
$arr[1,4] = 'LP1912'
$ws.Range('A50:E51').Value = $arr

# ---- Sheet: 6203-6173 ----
$ws = $wb.Worksheets.Item('6203-6173')
$ws.Range('A2').Value = 'Última actualización: 20:12:55'
$ws.Range('A3').Value = 'Total filas: 59'

$arr = New-Object 'object[,]' 5,5
$arr[0,0] = '20:12:55'
$arr[0,1] = '20:38'
$arr[0,2] = '215A_LA PLATA'
$arr[0,3] = 26
$arr[0,4] = 'L6173'
$arr[1,0] = '18:49:07'
$arr[1,1] = '20:39'
$arr[1,2] = '215A_LA PLATA'
$arr[1,3] = 110
$arr[1,4] = 'L6173'
$arr[2,0] = '19:52:49'
$arr[2,1] = '21:28'
$arr[2,2] = '215C_LA PLATA'
$arr[2,3] = 96
$arr[2,4] = 'L6203'
$arr[3,0] = '19:39:27'
$arr[3,1] = '21:29'
$arr[3,2] = '215C_LA PLATA'
$arr[3,3] = 110
$arr[3,4] = 'L6203'
$arr[4,0] = '20:12:55'
$arr[4,1] = '22:04'
$arr[4,2] = '215A_LA PLATA'
$arr[4,3] = 112
$arr[4,4] = 'L6173'
$ws.Range('A60:E64').Value = $arr

